# Commit: "update the data, lack dali and baselie, rerun dali"
#
# The baseline (B) and DALI (C) columns are cleared out because that data
# is being re-run / is currently missing, and the INR column (D) values
# are updated from 559 to 551 for every data row (rows 2-9).
# The current selection is also left on the cleared B2:C9 range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update INR (column D) values for rows 2-9 to the new rerun figure.
$ws.Range("D2:D9").Value = 551

# Clear out baseline (B) and DALI (C) columns - data is missing / being rerun.
$ws.Range("B2:C9").ClearContents()

# Leave the selection on the cleared range, matching the active cell state.
$ws.Range("B2:C9").Select()
